$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 57
$ws.Range("D2").Value = 50
$ws.Range("F2").Value = 0.0005
$ws.Range("H2").Value = 10

$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"
$ws.Range("K2").Value = 8
$ws.Range("L2").Value = 2

$ws.Range("D2").Select()
